$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("week")

$ws.Range("S2").Value = 301
$ws.Range("T2").Value = 380
$ws.Range("U2").Value = 426
$ws.Range("V2").Value = 328
$ws.Range("W2").Value = 321
$ws.Range("X2").Value = 440
$ws.Range("Y2").Value = 310
$ws.Range("Z2").Value = 323
$ws.Range("AA2").Value = 404
$ws.Range("AB2").Value = 402
$ws.Range("AC2").Value = 422
$ws.Range("AD2").Value = 336
$ws.Range("AE2").Value = 412
$ws.Range("AF2").Value = 328
$ws.Range("AG2").Value = 385
